# Apply the updated Gemini API keys to column D (gemini_api_key) of Sheet1
# and adjust the view/selection + ensure a new style entry for D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New individual API keys for rows 2-5 and 7-9 (values previously shared one of
# two keys, now each row gets a distinct key)
$ws.Range("D2").Value  = "AIzaSyBjQ36QX6K6cDCCGhkTu2ClKEiNk8frhzc"
$ws.Range("D3").Value  = "AIzaSyCKwFG2UGov5YAV0Hqfsd2vTy44KFttR74"
$ws.Range("D4").Value  = "AIzaSyAedzUac3wK1D-gU7QTAijLN16aMUKFGHU"
$ws.Range("D5").Value  = "AIzaSyDqNPQAdZ_cC9Y6ZvHaZhm-DUkTDejS78g"
$ws.Range("D7").Value  = "AIzaSyAd_W-jZ6sloI91NNB80gUj6fbHSy5n9gw"
$ws.Range("D8").Value  = "AIzaSyDGAls72zYru-Y9qa9J8vLREgO27J0N5kk"
$ws.Range("D9").Value  = "AIzaSyD9ElFE_jt973Wc-PjuZtLEWLBUTULYl70"

# D10 keeps the original key value, but gets a new font style applied
$ws.Range("D10").Value = "AIzaSyAug2O8re4wQ8LkM2LWqL8k792y5biZLHA"
$ws.Range("D10").Font.Name = "Calibri"

# Update the view: scroll position and current selection
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
